$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values (runs, balls, fours, sixes) per row, taken from the diff's
# "after" state. Row 6 is left untouched (no change in the diff).
$targets = @{
    2 = @("33", "16", "2", "2")
    3 = @("23", "15", "2", "1")
    4 = @("8",  "7",  "1", "0")
    5 = @("12", "7",  "1", "1")
    7 = @("0",  "0",  "0", "0")
    8 = @("7",  "5",  "1", "0")
    9 = @("20", "9",  "1", "2")
}

$cols = @("C", "D", "E", "F")

foreach ($row in $targets.Keys) {
    $vals = $targets[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "$($cols[$i])$row"
        $cell = $ws.Range($addr)
        # Keep these as text cells (matching the workbook's existing
        # "number stored as text" convention) rather than letting Excel
        # coerce the numeric-looking string into a true number.
        $cell.NumberFormat = "@"
        $cell.Value = $vals[$i]
        $cell.Style = "Normal"
    }
}
